$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: add new note about prefix (this string should become shared-string index 23)
$ws.Range("B11").Value = "Dont invoke if the message starts with the Guild's prefix"

# Row 3: consolidate B3/C3 into B3 (this string should become shared-string index 24), remove C3
$ws.Range("C3").ClearContents()
$ws.Range("B3").Value = "Add multiple Users as parameter, high5, not more then 1 mention of same person"
